# Update the "Avverkningsanmälningar" sheet:
#  - bump the "Förändrad" date (column C) from 2023-09-11 (45182) to
#    2023-09-13 (45184) for every existing data row (2..153)
#  - give row 153 an explicit row height (matches newly-added rows below)
#  - append four new data rows (154..157) for cases A 42802/42803/42801/42804-2023

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. bulk-update the "changed" date column for all existing rows ---
$ws.Range("C2:C153").Value2 = 45184

# --- 2. row 153 now carries an explicit (default) row height ---
$ws.Rows.Item(153).RowHeight = 15

# --- 3. append the four new rows ---
$newRows = @(
    @{ Row = 154; Beteckning = "A 42802-2023"; Area = 0.6; CustomHeight = $true },
    @{ Row = 155; Beteckning = "A 42803-2023"; Area = 0.7; CustomHeight = $true },
    @{ Row = 156; Beteckning = "A 42801-2023"; Area = 0.8; CustomHeight = $true },
    @{ Row = 157; Beteckning = "A 42804-2023"; Area = 0.6; CustomHeight = $false }
)

foreach ($entry in $newRows) {
    $r = $entry.Row

    $ws.Cells.Item($r, 1).Value2 = $entry.Beteckning                 # A: Beteckning

    $ws.Cells.Item($r, 2).Value2 = 45182                              # B: Datum
    $ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 3).Value2 = 45184                              # C: Förändrad
    $ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"

    $ws.Cells.Item($r, 4).Value2 = "ÖREBRO LÄN"                       # D: Län
    $ws.Cells.Item($r, 5).Value2 = "KARLSKOGA"                        # E: Kommun
    $ws.Cells.Item($r, 6).Value2 = "Sveaskog"                         # F: Markägare
    $ws.Cells.Item($r, 7).Value2 = $entry.Area                        # G: Area (ha)

    for ($c = 8; $c -le 17; $c++) {                                   # H..Q: counts, all 0
        $ws.Cells.Item($r, $c).Value2 = 0
    }

    $ws.Cells.Item($r, 18).WrapText = $true                           # R: Artnamn (wrap style)

    if ($entry.CustomHeight) {
        $ws.Rows.Item($r).RowHeight = 15
    }
}
